$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells hold values like "27.157.94" or "1.0000" that Excel
# would otherwise re-interpret as numbers (losing digits/trailing zeros), so we
# force those specific cells to Text format right before writing their value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.157.94'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.904.38'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.01'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5230'
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3761'
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07248'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.14'
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9025'
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08549'
$ws.Range("E12").Value = '  +11.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.902.14'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '95.19'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.290'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008635'
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.56'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.193.35'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.069'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.152.44'
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.429'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.289'
$ws.Range("E25").Value = '  +3.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.11'
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.751'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.22'
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.97'
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.812'
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.899'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09259'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("E33").Value = '  +4.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05048'
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.236'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.456'
$ws.Range("E36").Value = '  +5.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.948'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.618'
$ws.Range("E38").Value = '  +0.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5713'
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01998'
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.079'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.639'
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '115.97'
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4862'
$ws.Range("E46").Value = '  +1.19%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9998'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.10'
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.613'
$ws.Range("E49").Value = '  +1.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.49'
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.09'
